# The deck ships with two embedded themes:
#   theme1.xml -> "Integral"      (used by the slide master / all slides)
#   theme2.xml -> "Office Theme"  (used by the notes master)
#
# The authored change swaps which colour scheme backs the slides: the
# slide master (theme1.xml) is re-coloured with the stock "Office Theme"
# palette. We reproduce that by pushing the 12 standard theme colours
# (Background/Text 1-2, Accent 1-6, Hyperlink, Followed Hyperlink) that
# make up the "Office Theme" colour scheme onto the presentation's
# ThemeColorScheme, which PowerPoint stores once on the shared slide
# master theme (so every slide picks it up, since they all follow the
# same master).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Order matches MsoThemeColorSchemeIndex 1..12:
# Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink.
# Values are VBA-style RGB() integers (0xBBGGRR) for the "Office Theme"
# palette (44546A, E7E6E6, 5B9BD5, ED7D31, A5A5A5, FFC000, 4472C4,
# 70AD47, 0563C1, 954F72 alongside the standard black/white).
$officeThemeColors = @(
    0,          # Dark1   000000
    16777215,   # Light1  FFFFFF
    6968388,    # Dark2   44546A
    15132391,   # Light2  E7E6E6
    13998939,   # Accent1 5B9BD5
    3243501,    # Accent2 ED7D31
    10855845,   # Accent3 A5A5A5
    49407,      # Accent4 FFC000
    12874308,   # Accent5 4472C4
    4697456,    # Accent6 70AD47
    12673797,   # Hyperlink         0563C1
    7491477     # Followed Hyperlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
